$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.675.57"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "1.883.19"
$ws.Range("E3").Value = "  -1.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -1.73%  "

$ws.Range("E5").Value = "  -1.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  -2.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5123"
$ws.Range("E7").Value = "  -1.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3958"
$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08443"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.118"
$ws.Range("E10").Value = "  -1.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.312"
$ws.Range("E11").Value = "  +0.28%  "

$ws.Range("D12").Value = "1.884.77"
$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.56"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.302"
$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("E15").Value = "  -1.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001109"
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.62"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06744"
$ws.Range("E18").Value = "  -1.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.80"
$ws.Range("E19").Value = "  -0.85%  "

$ws.Range("E20").Value = "  -1.86%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.989"
$ws.Range("E21").Value = "  -1.64%  "

$ws.Range("D22").Value = "28.706.72"
$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.17"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.255"
$ws.Range("E24").Value = "  -0.25%  "

$ws.Range("D25").Value = "2.101.89"
$ws.Range("E25").Value = "  -0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.36"
$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.78"
$ws.Range("E27").Value = "  -0.93%  "

$ws.Range("E28").Value = "  -2.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "127.26"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("E30").Value = "  -1.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.053"
$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.822"
$ws.Range("E32").Value = "  -2.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.629"
$ws.Range("E33").Value = "  -1.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.02473"
$ws.Range("E34").Value = "  +0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06548"
$ws.Range("E35").Value = "  -1.57%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2200"
$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.972"
$ws.Range("E37").Value = "  -4.96%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.266"
$ws.Range("E38").Value = "  +0.41%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.203"
$ws.Range("E39").Value = "  +0.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.115"
$ws.Range("E40").Value = "  +2.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6479"
$ws.Range("E41").Value = "  -1.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.22"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.009"
$ws.Range("E43").Value = "  -1.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6090"
$ws.Range("E44").Value = "  -0.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.06"
$ws.Range("E45").Value = "  -1.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.706"
$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.036"
$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.224"
$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.78"
$ws.Range("E49").Value = "  -0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.196"
$ws.Range("E50").Value = "  -7.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06863"
$ws.Range("E51").Value = "  -1.37%  "
